$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 17 de Octubre de 2020 a las 21:20"

# Update country data rows (refreshed values + three country reorderings)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 8324898
$ws.Cells.Item(4, 3).Value = 36620
$ws.Cells.Item(4, 4).Value = 5412256
$ws.Cells.Item(4, 5).Value = 2688558
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 440
$ws.Cells.Item(4, 8).Value = 224084

$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 7492548
$ws.Cells.Item(5, 3).Value = 61913
$ws.Cells.Item(5, 4).Value = 6594155
$ws.Cells.Item(5, 5).Value = 784329
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1032
$ws.Cells.Item(5, 8).Value = 114064

$ws.Cells.Item(11, 1).Value = "Francia"
$ws.Cells.Item(11, 2).Value = 867197
$ws.Cells.Item(11, 3).Value = 32427
$ws.Cells.Item(11, 4).Value = 104696
$ws.Cells.Item(11, 5).Value = 729109
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 89
$ws.Cells.Item(11, 8).Value = 33392

$ws.Cells.Item(12, 1).Value = "Peru"
$ws.Cells.Item(12, 2).Value = 862417
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 769077
$ws.Cells.Item(12, 5).Value = 59692
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 33648

$ws.Cells.Item(13, 1).Value = "Mexico"
$ws.Cells.Item(13, 2).Value = 841661
$ws.Cells.Item(13, 3).Value = 6751
$ws.Cells.Item(13, 4).Value = 612216
$ws.Cells.Item(13, 5).Value = 143741
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 419
$ws.Cells.Item(13, 8).Value = 85704

$ws.Cells.Item(21, 1).Value = "Alemania"
$ws.Cells.Item(21, 2).Value = 361172
$ws.Cells.Item(21, 3).Value = 4380
$ws.Cells.Item(21, 4).Value = 290000
$ws.Cells.Item(21, 5).Value = 61323
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 13
$ws.Cells.Item(21, 8).Value = 9849

$ws.Cells.Item(52, 1).Value = "Etiopia"
$ws.Cells.Item(52, 2).Value = 88434
$ws.Cells.Item(52, 3).Value = 600
$ws.Cells.Item(52, 4).Value = 42099
$ws.Cells.Item(52, 5).Value = 44989
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 9
$ws.Cells.Item(52, 8).Value = 1346

$ws.Cells.Item(90, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(90, 2).Value = 23201
$ws.Cells.Item(90, 3).Value = 594
$ws.Cells.Item(90, 4).Value = 17111
$ws.Cells.Item(90, 5).Value = 5266
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 3
$ws.Cells.Item(90, 8).Value = 824

$ws.Cells.Item(92, 1).Value = "Costa de Marfil"
$ws.Cells.Item(92, 2).Value = 20301
$ws.Cells.Item(92, 3).Value = 26
$ws.Cells.Item(92, 4).Value = 19983
$ws.Cells.Item(92, 5).Value = 197
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 121

$ws.Cells.Item(98, 1).Value = "Zambia"
$ws.Cells.Item(98, 2).Value = 15789
$ws.Cells.Item(98, 3).Value = 77
$ws.Cells.Item(98, 4).Value = 14927
$ws.Cells.Item(98, 5).Value = 516
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 346

$ws.Cells.Item(104, 1).Value = "Namibia"
$ws.Cells.Item(104, 2).Value = 12263
$ws.Cells.Item(104, 3).Value = 48
$ws.Cells.Item(104, 4).Value = 10419
$ws.Cells.Item(104, 5).Value = 1713
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 131

$ws.Cells.Item(106, 1).Value = "Maldivas"
$ws.Cells.Item(106, 2).Value = 11178
$ws.Cells.Item(106, 3).Value = 24
$ws.Cells.Item(106, 4).Value = 10097
$ws.Cells.Item(106, 5).Value = 1045
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = 36

$ws.Cells.Item(116, 1).Value = "Zimbabue"
$ws.Cells.Item(116, 2).Value = 8110
$ws.Cells.Item(116, 3).Value = 11
$ws.Cells.Item(116, 4).Value = 7673
$ws.Cells.Item(116, 5).Value = 206
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 231

$ws.Cells.Item(136, 1).Value = "Ruanda"
$ws.Cells.Item(136, 2).Value = 4971
$ws.Cells.Item(136, 3).Value = 6
$ws.Cells.Item(136, 4).Value = 4768
$ws.Cells.Item(136, 5).Value = 169
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 34

$ws.Cells.Item(140, 1).Value = "Aruba"
$ws.Cells.Item(140, 2).Value = 4304
$ws.Cells.Item(140, 3).Value = 15
$ws.Cells.Item(140, 4).Value = 3983
$ws.Cells.Item(140, 5).Value = 287
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 2
$ws.Cells.Item(140, 8).Value = 34

$ws.Cells.Item(150, 1).Value = "Mali"
$ws.Cells.Item(150, 2).Value = 3379
$ws.Cells.Item(150, 3).Value = 1
$ws.Cells.Item(150, 4).Value = 2570
$ws.Cells.Item(150, 5).Value = 677
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 132

$ws.Cells.Item(154, 1).Value = "Republica de Chipre"
$ws.Cells.Item(154, 2).Value = 2581
$ws.Cells.Item(154, 3).Value = 202
$ws.Cells.Item(154, 4).Value = 1444
$ws.Cells.Item(154, 5).Value = 1112
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 25

$ws.Cells.Item(155, 1).Value = "Benin"
$ws.Cells.Item(155, 2).Value = 2496
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 2330
$ws.Cells.Item(155, 5).Value = 125
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 41

$ws.Cells.Item(156, 1).Value = "Uruguay"
$ws.Cells.Item(156, 2).Value = 2450
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 2042
$ws.Cells.Item(156, 5).Value = 357
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 51

$ws.Cells.Item(157, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(157, 2).Value = 2389
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 1782
$ws.Cells.Item(157, 5).Value = 566
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 41

$ws.Cells.Item(166, 1).Value = "Republica del Chad"
$ws.Cells.Item(166, 2).Value = 1365
$ws.Cells.Item(166, 3).Value = 4
$ws.Cells.Item(166, 4).Value = 1179
$ws.Cells.Item(166, 5).Value = 93
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 93

$ws.Cells.Item(194, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(194, 2).Value = 150
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 120
$ws.Cells.Item(194, 5).Value = 27
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 1
$ws.Cells.Item(194, 8).Value = 3

$ws.Cells.Item(200, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(200, 2).Value = 67
$ws.Cells.Item(200, 3).Value = 2
$ws.Cells.Item(200, 4).Value = 64
$ws.Cells.Item(200, 5).Value = 3
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0

$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 1

$ws.Cells.Item(217, 1).Value = "Islas Malvinas"
$ws.Cells.Item(217, 2).Value = 13
$ws.Cells.Item(217, 3).Value = 0
$ws.Cells.Item(217, 4).Value = 13
$ws.Cells.Item(217, 5).Value = 0
$ws.Cells.Item(217, 6).Value = 0
$ws.Cells.Item(217, 7).Value = 0
$ws.Cells.Item(217, 8).Value = 0

